{"js": "// The title paragraph reads \"Version 2.\" and must go back to \"Version 1.\"\n// (re-reverting the \"Wireframes version 2\" bump). The original markup has\n// \"Versi\"/\"on\" split across two runs (spell-check artifact) followed by a\n// \" 2\" run, the _GoBack bookmark, and a trailing \".\" run. The target state\n// merges \"Versi\"+\"on\" into a single \"Version\" run, folds the new \"1.\" into\n// the \" 2\" run (now \" 1.\"), and removes the now-redundant trailing \".\" run.\nconst body = context.document.body;\n\n// 1) Merge \"Versi\" + \"on\" into a single \"Version\" run.\nconst versionHits = body.search(\"Version\", { matchCase: true });\nversionHits.load(\"text\");\nawait context.sync();\nif (versionHits.items.length === 0) {\n  throw new Error(\"Could not find 'Version' in the document body.\");\n}\nversionHits.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n\n// 2) Replace the standalone \"2\" with \"1.\".\nconst twoHits = body.search(\"2\", { matchCase: true });\ntwoHits.load(\"text\");\nawait context.sync();\nif (twoHits.items.length === 0) {\n  throw new Error(\"Could not find the '2' to replace.\");\n}\ntwoHits.items[0].insertText(\"1.\", \"Replace\");\nawait context.sync();\n\n// 3) Delete the now-orphaned trailing \".\" run (the paragraph now reads\n//    \"Version 1..\" until this run is removed).\nconst dotHits = body.search(\".\", { matchCase: true });\ndotHits.load(\"text\");\nawait context.sync();\nif (dotHits.items.length === 0) {\n  throw new Error(\"Could not find the trailing '.' to remove.\");\n}\ndotHits.items[dotHits.items.length - 1].delete();\nawait context.sync();\n", "ps1": "# The title paragraph reads \"Version 2.\" and must go back to \"Version 1.\"\n# (re-reverting the \"Wireframes version 2\" bump). The original markup has\n# \"Versi\"/\"on\" split across two runs (spell-check artifact) followed by a\n# \" 2\" run, the _GoBack bookmark, and a trailing \".\" run. The target state\n# merges \"Versi\"+\"on\" into a single \"Version\" run, folds the new \"1.\" into\n# the \" 2\" run (now \" 1.\"), and removes the now-redundant trailing \".\" run\n# while keeping the _GoBack bookmark intact.\n$d = $word.ActiveDocument\n\n# 1) Merge \"Versi\" + \"on\" into a single \"Version\" run.\n$rngVersion = $d.Content\n$rngVersion.Find.Execute(\"Version\", $true, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2) | Out-Null\n\n# 2) Replace the standalone \"2\" with \"1.\".\n$rngTwo = $d.Content\n$rngTwo.Find.Execute(\"2\", $true, $false, $false, $false, $false, $true, 1, $false, \"1.\", 2) | Out-Null\n\n# 3) Remove the now-orphaned trailing \".\" run that follows the _GoBack\n#    bookmark, without disturbing the bookmark itself.\n$bm = $d.Bookmarks(\"_GoBack\")\n$para = $d.Paragraphs(1)\n$paraEnd = $para.Range.End\n$tail = $d.Range($bm.End, $paraEnd - 1)\nif ($tail.Text.Length -gt 0) {\n  $tail.Delete()\n}\n"}
